$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) load_shedding sheet: point the time-series columns at the "current"
#    (no-"future_"-prefix) data files instead of the future_* variants.
# ---------------------------------------------------------------------------
$wsLS = $wb.Worksheets.Item("load_shedding")
$wsLS.Range("D2").Value = "amiris-config/data/load.csv"
$wsLS.Range("D3").Value = "amiris-config/data/LS_high.csv"
$wsLS.Range("D5").Value = "amiris-config/data/LS_low.csv"
$wsLS.Range("D6").Value = "amiris-config/data/LS_mid.csv"

# ---------------------------------------------------------------------------
# 2) times sheet: shift Start/Stop time back (stochastic run now starts /
#    stops ~1461 "minutes"-scaled units earlier).
# ---------------------------------------------------------------------------
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 54788.99861111111
$wsTimes.Range("B3").Value = 55152.99861111111

# ---------------------------------------------------------------------------
# 3) conventionals sheet: a new OTHER unit (id 8 / 20243300061) is inserted
#    ahead of the existing units, pushing the rest down by one row.
# ---------------------------------------------------------------------------
$wsConv = $wb.Worksheets.Item("conventionals")

# Shift existing rows 7..9 down to 8..10 (copy bottom-up so data is never
# clobbered before it's been read).
$wsConv.Cells.Item(10,1).Value = $wsConv.Cells.Item(9,1).Value2
$wsConv.Cells.Item(10,2).Value = $wsConv.Cells.Item(9,2).Value2
$wsConv.Cells.Item(10,3).Value = $wsConv.Cells.Item(9,3).Value2
$wsConv.Cells.Item(10,4).Value = $wsConv.Cells.Item(9,4).Value2
$wsConv.Cells.Item(10,5).Value = $wsConv.Cells.Item(9,5).Value2
$wsConv.Cells.Item(10,6).Value = $wsConv.Cells.Item(9,6).Value2
$wsConv.Cells.Item(10,7).Value = $wsConv.Cells.Item(9,7).Value2
$wsConv.Cells.Item(10,1).Font.Bold = $true
$wsConv.Cells.Item(10,1).HorizontalAlignment = -4108
$wsConv.Cells.Item(10,1).VerticalAlignment = -4160
$wsConv.Cells.Item(10,1).Borders.LineStyle = 1

$wsConv.Cells.Item(9,1).Value = $wsConv.Cells.Item(8,1).Value2
$wsConv.Cells.Item(9,2).Value = $wsConv.Cells.Item(8,2).Value2
$wsConv.Cells.Item(9,3).Value = $wsConv.Cells.Item(8,3).Value2
$wsConv.Cells.Item(9,4).Value = $wsConv.Cells.Item(8,4).Value2
$wsConv.Cells.Item(9,5).Value = $wsConv.Cells.Item(8,5).Value2
$wsConv.Cells.Item(9,6).Value = $wsConv.Cells.Item(8,6).Value2
$wsConv.Cells.Item(9,7).Value = $wsConv.Cells.Item(8,7).Value2

$wsConv.Cells.Item(8,1).Value = $wsConv.Cells.Item(7,1).Value2
$wsConv.Cells.Item(8,2).Value = $wsConv.Cells.Item(7,2).Value2
$wsConv.Cells.Item(8,3).Value = $wsConv.Cells.Item(7,3).Value2
$wsConv.Cells.Item(8,4).Value = $wsConv.Cells.Item(7,4).Value2
$wsConv.Cells.Item(8,5).Value = $wsConv.Cells.Item(7,5).Value2
$wsConv.Cells.Item(8,6).Value = $wsConv.Cells.Item(7,6).Value2
$wsConv.Cells.Item(8,7).Value = $wsConv.Cells.Item(7,7).Value2

# Write the new unit into (now-vacated) row 7.
$wsConv.Cells.Item(7,1).Value = 8
$wsConv.Cells.Item(7,2).Value = 20243300061
$wsConv.Cells.Item(7,3).Value = "OTHER"
$wsConv.Cells.Item(7,4).Value = 1.5
$wsConv.Cells.Item(7,5).Value = 0.4
$wsConv.Cells.Item(7,6).Value = 1500
$wsConv.Cells.Item(7,7).Value = 1500

# ---------------------------------------------------------------------------
# 4) renewables sheet: one new WindOff-replacement OtherPV unit is inserted
#    ahead of id 34/35, and two more OtherPV units are appended after them
#    (ids 36 and 37), pushing the former last row (id 35 / WindOff) down to
#    become id 38.
# ---------------------------------------------------------------------------
$wsRen = $wb.Worksheets.Item("renewables")

# Move old row 37 (id 35, WindOff) down to row 40.
for ($col = 1; $col -le 9; $col++) {
    $wsRen.Cells.Item(40,$col).Value = $wsRen.Cells.Item(37,$col).Value2
}
$wsRen.Cells.Item(40,1).Font.Bold = $true
$wsRen.Cells.Item(40,1).HorizontalAlignment = -4108
$wsRen.Cells.Item(40,1).VerticalAlignment = -4160
$wsRen.Cells.Item(40,1).Borders.LineStyle = 1

# Move old row 36 (id 34, WindOff) down to row 37.
for ($col = 1; $col -le 9; $col++) {
    $wsRen.Cells.Item(37,$col).Value = $wsRen.Cells.Item(36,$col).Value2
}

# New row 36: id 34, OtherPV replacement unit ahead of the old id-34 unit.
$wsRen.Cells.Item(36,1).Value = 34
$wsRen.Cells.Item(36,2).Value = 20290300029
$wsRen.Cells.Item(36,3).Value = 3000
$wsRen.Cells.Item(36,4).Value = 0
$wsRen.Cells.Item(36,5).Value = "OtherPV"
$wsRen.Cells.Item(36,6).Value = "NONE"
$wsRen.Cells.Item(36,7).Value = "-"
$wsRen.Cells.Item(36,8).Value = "-"
$wsRen.Cells.Item(36,9).Value = "-"

# New row 38: id 36, OtherPV unit.
$wsRen.Cells.Item(38,1).Value = 36
$wsRen.Cells.Item(38,2).Value = 20270300028
$wsRen.Cells.Item(38,3).Value = 2500
$wsRen.Cells.Item(38,4).Value = 0
$wsRen.Cells.Item(38,5).Value = "OtherPV"
$wsRen.Cells.Item(38,6).Value = "NONE"
$wsRen.Cells.Item(38,7).Value = "-"
$wsRen.Cells.Item(38,8).Value = "-"
$wsRen.Cells.Item(38,9).Value = "-"
$wsRen.Cells.Item(38,1).Font.Bold = $true
$wsRen.Cells.Item(38,1).HorizontalAlignment = -4108
$wsRen.Cells.Item(38,1).VerticalAlignment = -4160
$wsRen.Cells.Item(38,1).Borders.LineStyle = 1

# New row 39: id 37, OtherPV unit.
$wsRen.Cells.Item(39,1).Value = 37
$wsRen.Cells.Item(39,2).Value = 20260300027
$wsRen.Cells.Item(39,3).Value = 1000
$wsRen.Cells.Item(39,4).Value = 0
$wsRen.Cells.Item(39,5).Value = "OtherPV"
$wsRen.Cells.Item(39,6).Value = "NONE"
$wsRen.Cells.Item(39,7).Value = "-"
$wsRen.Cells.Item(39,8).Value = "-"
$wsRen.Cells.Item(39,9).Value = "-"
$wsRen.Cells.Item(39,1).Font.Bold = $true
$wsRen.Cells.Item(39,1).HorizontalAlignment = -4108
$wsRen.Cells.Item(39,1).VerticalAlignment = -4160
$wsRen.Cells.Item(39,1).Borders.LineStyle = 1
